# Apply the GA timetable conflict-resolution update to every weekly sheet
# (Phong_Tuan_1 .. Phong_Tuan_15). All sheets receive the identical edit to
# rows 8-11 of the "Lịch sử dụng phòng" table.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$newD8 = "Lớp: CL05" + [char]10 + "Môn: Kỹ năng mềm" + [char]10 + "GV: Hoàng Thị E" + [char]10 + "(Lý thuyết)"
$newH9 = "Lớp: CL05" + [char]10 + "Môn: Tiếng Anh chuyên ngành" + [char]10 + "GV: Hoàng Thị E" + [char]10 + "(Lý thuyết)"
$newG10 = "Lớp: CL10" + [char]10 + "Môn: Tiếng Anh chuyên ngành" + [char]10 + "GV: Hoàng Thị E" + [char]10 + "(Lý thuyết)"
$newG11 = "Lớp: CL10" + [char]10 + "Môn: Kỹ năng mềm" + [char]10 + "GV: Hoàng Thị E" + [char]10 + "(Lý thuyết)"
$newA10 = "C2" + [char]10 + "(15:00-17:00)"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- Row 8 ---
    # B8: room changes R102 -> R103
    $ws.Range("B8").Value2 = "R103"

    # G8 currently holds the CL10/English class block (style 8, filled).
    # Copy its formatting onto D8, then fill D8 with the "Kỹ năng mềm" block
    # that used to live in F10, and clear G8 back to an empty (style 7) cell.
    $ws.Range("G8").Copy()
    $ws.Range("D8").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $ws.Range("D8").Value2 = $newD8

    $ws.Range("C8").Copy()
    $ws.Range("G8").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $ws.Range("G8").Value2 = ""

    # --- Row 9 ---
    # H9: teacher changes Võ Văn F -> Hoàng Thị E (class/subject unchanged)
    $ws.Range("H9").Value2 = $newH9
    # Row 9 has no explicit custom height in the source file (it relies on the
    # default row height); writing to it makes the engine stamp an explicit
    # height, so auto-fit it back down to a plain, non-custom row height.
    $ws.Rows.Item(9).AutoFit()

    # --- Row 10 ---
    # A10: slot changes C1 (13:00-15:00) -> C2 (15:00-17:00)
    $ws.Range("A10").Value2 = $newA10

    # F10 currently holds the CL05/Kỹ năng mềm block (style 8, filled).
    # Copy its formatting onto G10, fill G10 with the CL10/English block
    # (now taught by Hoàng Thị E), and clear F10 back to empty (style 7).
    $ws.Range("F10").Copy()
    $ws.Range("G10").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $ws.Range("G10").Value2 = $newG10

    $ws.Range("E10").Copy()
    $ws.Range("F10").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $ws.Range("F10").Value2 = ""

    # --- Row 11 ---
    # B11: room changes R104 -> R101
    $ws.Range("B11").Value2 = "R101"

    # G11: teacher changes Ngô Văn I -> Hoàng Thị E (class/subject unchanged)
    $ws.Range("G11").Value2 = $newG11
}
